$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189, shifting existing rows 189-216 down to 190-217.
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A189").Value2 = 10
$ws.Range("B189").Value2 = "Vega Modelo de Temuco"
$ws.Range("C189").Value2 = "La Araucanía"
$ws.Range("D189").Value2 = 44474
$ws.Range("E189").Value2 = 9
$ws.Range("F189").Value2 = "Fruta"
$ws.Range("G189").Value2 = 100108
$ws.Range("H189").Value2 = "Tropicales y subtropicales"
$ws.Range("I189").Value2 = 100108002
$ws.Range("J189").Value2 = "Mango"
$ws.Range("K189").Value2 = "Sin especificar"
$ws.Range("L189").Value2 = "Primera"
$ws.Range("M189").Value2 = 300
$ws.Range("N189").Value2 = 8000
$ws.Range("O189").Value2 = 8000
$ws.Range("P189").Value2 = 8000
$ws.Range("Q189").Value2 = "$/bandeja 4 kilos"
$ws.Range("R189").Value2 = "Perú"
$ws.Range("S189").Value2 = 2000
$ws.Range("T189").Value2 = 4
